# "update tasks for Tet Holiday"
# Set Owner = "Hung" and Status = "on processing" for rows 22, 27, 28 on Sheet1,
# copying formatting from an existing "Hung"-owner cell (C4) and an existing
# "on processing"-status cell (G23) so the new cells pick up the same
# visual style already used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ownerRows = @(22, 27, 28)
foreach ($r in $ownerRows) {
    $ws.Range("C4").Copy()
    $ws.Range("C" + $r).PasteSpecial(-4122)
    $ws.Range("C" + $r).Value = "Hung"

    $ws.Range("G23").Copy()
    $ws.Range("G" + $r).PasteSpecial(-4122)
    $ws.Range("G" + $r).Value = "on processing"
}

# Restore the scroll position / active selection recorded in the sheet view
# (user had scrolled further down and was working around row 35 afterwards).
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C35").Select()
